$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.477.08"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "'1.575.31"
$ws.Range("E3").Value = "  +0.77%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'1.000"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").Value = "'288.03"
$ws.Range("E6").Value = "  +0.47%  "

$ws.Range("D7").Value = "'0.3687"
$ws.Range("E7").Value = "  +0.89%  "

$ws.Range("D9").Value = "'0.3329"
$ws.Range("E9").Value = "  -0.66%  "

$ws.Range("D10").Value = "'1.150"
$ws.Range("E10").Value = "  +2.07%  "

$ws.Range("D11").Value = "'0.07559"
$ws.Range("E11").Value = "  +2.15%  "

$ws.Range("E12").Value = "  -0.08%  "

$ws.Range("D13").Value = "'20.78"
$ws.Range("E13").Value = "  -0.21%  "

$ws.Range("D14").Value = "'5.947"
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("D15").Value = "'6.959"

$ws.Range("D16").Value = "'1.569.31"
$ws.Range("E16").Value = "  +0.42%  "

$ws.Range("D17").Value = "'0.00001121"
$ws.Range("E17").Value = "  +2.07%  "

$ws.Range("D18").Value = "'88.02"
$ws.Range("E18").Value = "  -1.16%  "

$ws.Range("D19").Value = "'0.06738"
$ws.Range("E19").Value = "  -0.34%  "

$ws.Range("D20").Value = "'1.0000"
$ws.Range("E20").Value = "  -0.14%  "

$ws.Range("D21").Value = "'6.393"
$ws.Range("E21").Value = "  +1.41%  "

$ws.Range("D22").Value = "'16.54"
$ws.Range("E22").Value = "  +2.95%  "

$ws.Range("D24").Value = "'22.473.39"
$ws.Range("E24").Value = "  +0.42%  "

$ws.Range("E25").Value = "  +0.77%  "

$ws.Range("E26").Value = "  +4.14%  "

$ws.Range("E27").Value = "  +0.43%  "

$ws.Range("D28").Value = "'19.69"
$ws.Range("E28").Value = "  +0.46%  "

$ws.Range("D29").Value = "'4.988"
$ws.Range("E29").Value = "  -0.27%  "

$ws.Range("E30").Value = "  +1.92%  "

$ws.Range("D31").Value = "'1.748.98"
$ws.Range("E31").Value = "  +0.62%  "

$ws.Range("D32").Value = "'1.092"
$ws.Range("E32").Value = "  +3.06%  "

$ws.Range("D33").Value = "'6.111"
$ws.Range("E33").Value = "  +0.36%  "

$ws.Range("D34").Value = "'1.995"
$ws.Range("E34").Value = "  +0.11%  "

$ws.Range("E35").Value = "  +2.82%  "

$ws.Range("D36").Value = "'0.08365"
$ws.Range("E36").Value = "  +1.28%  "

$ws.Range("D37").Value = "'0.02467"

$ws.Range("E38").Value = "  +0.97%  "

$ws.Range("D39").Value = "'0.06396"
$ws.Range("E39").Value = "  +0.24%  "

$ws.Range("D40").Value = "'1.293"
$ws.Range("E40").Value = "  -1.28%  "

$ws.Range("D41").Value = "'5.368"
$ws.Range("E41").Value = "  +0.91%  "

$ws.Range("D42").Value = "'11.47"
$ws.Range("E42").Value = "  +2.98%  "

$ws.Range("D43").Value = "'0.6274"
$ws.Range("E43").Value = "  +3.43%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'14.06"
$ws.Range("E44").Value = "  +2.87%  "

$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("D46").Value = "'0.6117"
$ws.Range("E46").Value = "  +6.67%  "

$ws.Range("D47").Value = "'3.780"
$ws.Range("E47").Value = "  +0.35%  "

$ws.Range("D48").Value = "'2.058"
$ws.Range("E48").Value = "  +2.28%  "

$ws.Range("D49").Value = "'125.38"
$ws.Range("E49").Value = "  +0.43%  "

$ws.Range("E50").Value = "  +0.02%  "

$ws.Range("D51").Value = "'0.07223"
$ws.Range("E51").Value = "  -0.08%  "
